$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The values in row 2 (FacturaID, FechaEmision, personeriaJuridica, CondicionVenta,
# PlazoCredito, MedioPago) are numeric/zero-padded looking strings that must stay
# plain text (leading zeros, >15-digit IDs). Force text format before assigning so
# Excel's COM layer doesn't auto-coerce them into numbers, then restore the
# original (unstyled/General) look so no visible formatting changes.
$dataRow = $ws.Range("A2:F2")
$dataRow.NumberFormat = "@"

$ws.Range("A2").Value = "00100102010000009206"
$ws.Range("B2").Value = "2025-03-31T10:24:00-06:00"
$ws.Range("C2").Value = "3101135332"
$ws.Range("D2").Value = "02"
$ws.Range("E2").Value = "15"
$ws.Range("F2").Value = "04"

$dataRow.ClearFormats()
